{"js": "// Adds explanatory \"//\" comments above several methods/blocks in the\n// Red.cpp listing that lives in this document (one code line per\n// paragraph), and tidies up the blank-paragraph run that precedes\n// `void Red::agregarLinea(...)`, plus adds one trailing blank paragraph\n// at the very end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Helper: find the Nth (0-based `occurrence`) paragraph whose text equals\n// `text` exactly.\nfunction findParagraph(text, occurrence = 0) {\n  let seen = 0;\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) {\n      if (seen === occurrence) return items[i];\n      seen++;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + text + \" (occurrence \" + occurrence + \")\");\n}\n\n// Each entry: comment text to insert, BEFORE the paragraph matching\n// `anchorText` (Nth occurrence, 0-based).\nconst insertions = [\n  { anchorText: \"    this->nombre = nombre;\", occurrence: 0,\n    comment: \"    // Constructor de la clase Red\" },\n  { anchorText: \"    liberarMemoria(cabezaLineas);\", occurrence: 0,\n    comment: \"    // Destructor de la clase Red\" },\n  { anchorText: \"    return cabezaLineas;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para obtener un puntero a la cabeza de la lista de l\u00edneas\" },\n  { anchorText: \"    return nombre;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para obtener el nombre de la red\" },\n  { anchorText: \"    return numLineas;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para obtener el n\u00famero de l\u00edneas en la red\" },\n  { anchorText: \"    if (linea == nullptr) {\", occurrence: 0,\n    comment: \"    // Funci\u00f3n privada para liberar la memoria de las l\u00edneas\" },\n  { anchorText: \"    if (cabezaLineas == nullptr) {\", occurrence: 0,\n    comment: \"    // M\u00e9todo para agregar una nueva l\u00ednea a la red\" },\n  { anchorText: \"    if (cabezaLineas == nullptr) {\", occurrence: 1,\n    comment: \"    // M\u00e9todo para eliminar una l\u00ednea de la red\" },\n  { anchorText: \"    int totalEstaciones = 0;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para contar el n\u00famero total de estaciones en la red\" },\n  { anchorText: \"    int tiempo = 0;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para calcular el tiempo entre dos estaciones\" },\n  { anchorText: \"    Estacion* inicioEstacion = nullptr;\", occurrence: 0,\n    comment: \"    // M\u00e9todo para calcular el tiempo de llegada entre dos estaciones\" },\n];\n\nfor (const { anchorText, occurrence, comment } of insertions) {\n  const anchor = findParagraph(anchorText, occurrence);\n  anchor.insertParagraph(comment, \"Before\");\n}\n\n// Collapse the 3 consecutive blank paragraphs that precede\n// `void Red::agregarLinea(Linea* linea) {` down to a single blank\n// paragraph (remove two of them).\nconst agregarLinea = findParagraph(\"void Red::agregarLinea(Linea* linea) {\", 0);\nlet prev = agregarLinea.getPrevious();\nprev.load(\"text\");\nawait context.sync();\nlet removed = 0;\nwhile (removed < 2 && prev.text === \"\") {\n  const toDelete = prev;\n  prev = prev.getPrevious();\n  prev.load(\"text\");\n  toDelete.delete();\n  await context.sync();\n  removed++;\n}\n\n// Append one new blank paragraph right after the final closing brace of\n// `calcularTiempoLlegada` (i.e. before the pre-existing trailing blank\n// paragraph at the end of the document).\nconst veryLast = items[items.length - 1];\nveryLast.insertParagraph(\"\", \"Before\");\n", "ps1": "# Adds explanatory \"//\" comments above several methods/blocks in the\n# Red.cpp listing that lives in this document (one code line per\n# paragraph), tidies up the blank-paragraph run that precedes\n# `void Red::agregarLinea(...)`, and adds one trailing blank paragraph\n# at the very end of the document.\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($Text, $Occurrence) {\n    $seen = 0\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $t = $d.Paragraphs.Item($i).Range.Text\n        $t = $t.TrimEnd([char]13, [char]7)\n        if ($t -eq $Text) {\n            if ($seen -eq $Occurrence) {\n                return $i\n            }\n            $seen++\n        }\n    }\n    throw \"Paragraph not found: $Text (occurrence $Occurrence)\"\n}\n\n# Each row: anchor paragraph text / occurrence (0-based) to insert the new\n# comment paragraph BEFORE.\n$insertions = @(\n    @(\"    this->nombre = nombre;\", 0, \"    // Constructor de la clase Red\"),\n    @(\"    liberarMemoria(cabezaLineas);\", 0, \"    // Destructor de la clase Red\"),\n    @(\"    return cabezaLineas;\", 0, \"    // M\u00e9todo para obtener un puntero a la cabeza de la lista de l\u00edneas\"),\n    @(\"    return nombre;\", 0, \"    // M\u00e9todo para obtener el nombre de la red\"),\n    @(\"    return numLineas;\", 0, \"    // M\u00e9todo para obtener el n\u00famero de l\u00edneas en la red\"),\n    @(\"    if (linea == nullptr) {\", 0, \"    // Funci\u00f3n privada para liberar la memoria de las l\u00edneas\"),\n    @(\"    if (cabezaLineas == nullptr) {\", 0, \"    // M\u00e9todo para agregar una nueva l\u00ednea a la red\"),\n    @(\"    if (cabezaLineas == nullptr) {\", 1, \"    // M\u00e9todo para eliminar una l\u00ednea de la red\"),\n    @(\"    int totalEstaciones = 0;\", 0, \"    // M\u00e9todo para contar el n\u00famero total de estaciones en la red\"),\n    @(\"    int tiempo = 0;\", 0, \"    // M\u00e9todo para calcular el tiempo entre dos estaciones\"),\n    @(\"    Estacion* inicioEstacion = nullptr;\", 0, \"    // M\u00e9todo para calcular el tiempo de llegada entre dos estaciones\")\n)\n\nforeach ($row in $insertions) {\n    $anchorText = $row[0]\n    $occurrence = $row[1]\n    $comment = $row[2]\n    $idx = Find-ParagraphIndex $anchorText $occurrence\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphBefore()\n    # The freshly inserted (still empty) paragraph now sits at $idx; the\n    # original anchor paragraph shifted to $idx + 1.\n    $d.Paragraphs.Item($idx).Range.Text = $comment\n}\n\n# Collapse the 3 consecutive blank paragraphs that precede\n# `void Red::agregarLinea(Linea* linea) {` down to a single blank\n# paragraph (remove two of them).\n$agregarIdx = Find-ParagraphIndex \"void Red::agregarLinea(Linea* linea) {\" 0\n$removed = 0\nwhile ($removed -lt 2) {\n    $prevIdx = $agregarIdx - 1\n    $prevText = $d.Paragraphs.Item($prevIdx).Range.Text.TrimEnd([char]13, [char]7)\n    if ($prevText -ne \"\") {\n        break\n    }\n    $d.Paragraphs.Item($prevIdx).Range.Delete()\n    $agregarIdx = $agregarIdx - 1\n    $removed++\n}\n\n# Append one new blank paragraph right after the final closing brace of\n# `calcularTiempoLlegada` (i.e. before the pre-existing trailing blank\n# paragraph at the end of the document).\n$lastIdx = $d.Paragraphs.Count\n$d.Paragraphs.Item($lastIdx).Range.InsertParagraphBefore()\n"}
